# Updated symbol list on Tue Jan 24 03:34:00 UTC 2023 with GitHub Actions
# Applies refreshed price / volume(1h) data to the cryptos sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data cells in columns D (Price) and E (Volume(1h)) are stored as text
# strings in the source workbook (e.g. "319.26", "4.44%"). Force the
# number format to Text before writing so Excel keeps the literal string
# instead of re-interpreting it as a number/percentage, then restore the
# original (default) style so formatting is left unchanged.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "321.50"
$ws.Range("E2").Value = "5.15%"
$ws.Range("D3").Value = "36.03"
$ws.Range("E3").Value = "-0.41%"
$ws.Range("D4").Value = "5.132"
$ws.Range("E4").Value = "1.85%"
$ws.Range("D5").Value = "0.08073"
$ws.Range("E5").Value = "3.22%"
$ws.Range("D6").Value = "2.161"
$ws.Range("E6").Value = "1.79%"
$ws.Range("E7").Value = "1.80%"
$ws.Range("D8").Value = "4.132"
$ws.Range("E8").Value = "1.51%"
$ws.Range("E9").Value = "1.19%"
$ws.Range("D10").Value = "0.1007"
$ws.Range("E10").Value = "5.12%"
$ws.Range("D11").Value = "0.1881"
$ws.Range("E11").Value = "0.57%"
$ws.Range("D12").Value = "0.09218"
$ws.Range("E12").Value = "6.09%"
$ws.Range("D13").Value = "0.03602"
$ws.Range("E13").Value = "2.99%"
$ws.Range("D14").Value = "0.09944"
$ws.Range("E14").Value = "0.27%"
$ws.Range("D15").Value = "0.001442"
$ws.Range("E15").Value = "0.83%"
$ws.Range("D16").Value = "0.005671"
$ws.Range("E16").Value = "-1.01%"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").Value = "0.06%"
$ws.Range("E18").Value = "18.41%"
$ws.Range("E19").Value = "-1.29%"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").Value = "-1.13%"
$ws.Range("D21").Value = "5.054"
$ws.Range("E21").Value = "6.11%"
$ws.Range("E22").Value = "-3.84%"
$ws.Range("D23").Value = "0.04601"
$ws.Range("E23").Value = "-0.16%"
$ws.Range("E24").Value = "1.13%"
$ws.Range("D25").Value = "0.004745"
$ws.Range("E25").Value = "-6.98%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "-7.14%"
$ws.Range("D27").Value = "0.0004505"
$ws.Range("E27").Value = "65.28%"
$ws.Range("D39").Value = "0.01955"
$ws.Range("E39").Value = "7.43%"
$ws.Range("D40").Value = "0.04975"
$ws.Range("E40").Value = "4.62%"
$ws.Range("D41").Value = "0.007816"
$ws.Range("E41").Value = "4.04%"
$ws.Range("E42").Value = "-0.08%"
$ws.Range("D43").Value = "0.007818"
$ws.Range("E43").Value = "1.37%"
$ws.Range("D44").Value = "0.002068"
$ws.Range("E44").Value = "-7.33%"
$ws.Range("D45").Value = "0.01175"
$ws.Range("E45").Value = "8.84%"
$ws.Range("D46").Value = "0.00006317"
$ws.Range("E46").Value = "1.23%"
$ws.Range("E47").Value = "0.00%"
$ws.Range("E48").Value = "10.15%"
$ws.Range("E49").Value = "-4.95%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "0.00%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "0.00%"

# Restore original cell style/formatting now that the text values are set.
$dataRange.Style = $origStyle
